$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for team record, matching the style of the existing
# header row (bold, bordered, centered) by copying formats from AA1.
$ws.Range("AA1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record (W/L/T) repeated for every player row.
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 30).Value = 79
    $ws.Cells.Item($r, 31).Value = 83
    $ws.Cells.Item($r, 32).Value = 0
}
